# Update the "Generate Report for Handback" timestamps.
# These cells hold plain text (not real date serials) that look like
# "yyyy-mm-dd HH:mm:ss", so we set them as strings to keep the shared-string
# text type intact.

$wb = $excel.ActiveWorkbook

# Overview sheet: row for 5af550af-...md -> "Latest HO Xliff Generate Date" (col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 12:49:18"

# zh-cn sheet: row for 5af550af-...md -> Correspond Handoff Datetime (col H) / Correspond Handback DateTime (col K)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-03 12:49:13"
$wsZhCn.Range("K3").Value = "2016-09-03 12:49:32"

# de-de sheet: row for 5af550af-...md -> Correspond Handback DateTime (col K)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-09-03 12:49:39"
